$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2..297 from 45192 to 45202
for ($r = 2; $r -le 297; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# Add new row 298 with data
$ws.Cells.Item(298, 1).Value = "A 45833-2023"
$ws.Cells.Item(298, 2).Value = 45195
$ws.Cells.Item(298, 3).Value = 45202
$ws.Cells.Item(298, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item(298, 5).Value = "VANSBRO"
$ws.Cells.Item(298, 7).Value = 7.2
$ws.Cells.Item(298, 8).Value = 0
$ws.Cells.Item(298, 9).Value = 0
$ws.Cells.Item(298, 10).Value = 0
$ws.Cells.Item(298, 11).Value = 0
$ws.Cells.Item(298, 12).Value = 0
$ws.Cells.Item(298, 13).Value = 0
$ws.Cells.Item(298, 14).Value = 0
$ws.Cells.Item(298, 15).Value = 0
$ws.Cells.Item(298, 16).Value = 0
$ws.Cells.Item(298, 17).Value = 0

$ws.Cells.Item(298, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(298, 3).NumberFormat = "YYYY-MM-DD"
$ws.Range("R298").WrapText = $true

$ws.Rows.Item(297).RowHeight = 15

